$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.928.59"
$ws.Range("E2").Value = "'  -0.88%  "
$ws.Range("D3").Value = "'1.834.03"
$ws.Range("E3").Value = "'  -1.28%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'245.23"
$ws.Range("E5").Value = "'  +1.32%  "
$ws.Range("D6").Value = "'0.6933"
$ws.Range("E6").Value = "'  -1.00%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("D8").Value = "'0.07684"
$ws.Range("E8").Value = "'  -1.59%  "
$ws.Range("D9").Value = "'0.3046"
$ws.Range("E9").Value = "'  -2.18%  "
$ws.Range("D10").Value = "'23.34"
$ws.Range("E10").Value = "'  -2.76%  "
$ws.Range("D11").Value = "'0.07809"
$ws.Range("E11").Value = "'  +0.06%  "
$ws.Range("D12").Value = "'93.15"
$ws.Range("E12").Value = "'  +1.09%  "
$ws.Range("D13").Value = "'1.832.51"
$ws.Range("E13").Value = "'  -1.22%  "
$ws.Range("D14").Value = "'5.103"
$ws.Range("E14").Value = "'  -0.52%  "
$ws.Range("D15").Value = "'0.6815"
$ws.Range("E15").Value = "'  -1.45%  "
$ws.Range("D16").Value = "'6.587"
$ws.Range("E16").Value = "'  +0.29%  "
$ws.Range("D17").Value = "'0.000008256"
$ws.Range("E17").Value = "'  -2.96%  "
$ws.Range("D18").Value = "'28.923.58"
$ws.Range("E18").Value = "'  -0.92%  "
$ws.Range("D19").Value = "'240.94"
$ws.Range("E19").Value = "'  -3.04%  "
$ws.Range("D20").Value = "'2.074.72"
$ws.Range("E20").Value = "'  -1.35%  "
$ws.Range("D21").Value = "'12.69"
$ws.Range("E21").Value = "'  -1.88%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("D23").Value = "'7.458"
$ws.Range("E23").Value = "'  -1.51%  "
$ws.Range("D24").Value = "'1.0000"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("D25").Value = "'0.1502"
$ws.Range("E25").Value = "'  -2.17%  "
$ws.Range("D26").Value = "'158.36"
$ws.Range("E26").Value = "'  -1.33%  "
$ws.Range("D27").Value = "'8.754"
$ws.Range("E27").Value = "'  -1.86%  "
$ws.Range("D28").Value = "'18.17"
$ws.Range("E28").Value = "'  -2.23%  "
$ws.Range("D29").Value = "'1.538"
$ws.Range("E29").Value = "'  -2.04%  "
$ws.Range("D30").Value = "'4.212"
$ws.Range("E30").Value = "'  -1.52%  "
$ws.Range("D31").Value = "'4.157"
$ws.Range("E31").Value = "'  -2.00%  "
$ws.Range("D32").Value = "'1.193"
$ws.Range("E32").Value = "'  -1.07%  "
$ws.Range("D33").Value = "'0.05109"
$ws.Range("E33").Value = "'  -2.58%  "
$ws.Range("D34").Value = "'0.7792"
$ws.Range("E34").Value = "'  +2.78%  "
$ws.Range("D35").Value = "'1.854"
$ws.Range("E35").Value = "'  -0.96%  "
$ws.Range("D36").Value = "'1.144"
$ws.Range("E36").Value = "'  -2.69%  "
$ws.Range("D37").Value = "'2.696"
$ws.Range("E37").Value = "'  -0.19%  "
$ws.Range("D38").Value = "'1.290.59"
$ws.Range("E38").Value = "'  +4.84%  "
$ws.Range("D39").Value = "'0.01858"
$ws.Range("E39").Value = "'  +0.28%  "
$ws.Range("D40").Value = "'2.699"
$ws.Range("E40").Value = "'  -1.32%  "
$ws.Range("D41").Value = "'0.9548"
$ws.Range("E41").Value = "'  +5.91%  "
$ws.Range("D42").Value = "'6.139"
$ws.Range("E42").Value = "'  +4.81%  "
$ws.Range("D43").Value = "'106.78"
$ws.Range("E43").Value = "'  -3.06%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("D45").Value = "'9.685"
$ws.Range("E45").Value = "'  +1.80%  "
$ws.Range("D46").Value = "'0.00000000123"
$ws.Range("E46").Value = "'  -0.88%  "
$ws.Range("D47").Value = "'0.5166"
$ws.Range("E47").Value = "'  -0.26%  "
$ws.Range("D48").Value = "'1.974.98"
$ws.Range("E48").Value = "'  -1.57%  "
$ws.Range("D49").Value = "'63.90"
$ws.Range("E49").Value = "'  -6.41%  "
$ws.Range("D50").Value = "'1.753"
$ws.Range("E50").Value = "'  -0.65%  "
$ws.Range("D51").Value = "'6.959"
$ws.Range("E51").Value = "'  -0.79%  "

$ws.Range("D2:E51").Style = "Normal"
